$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New scores for column C (rows 1-8). These must land in the workbook as
# plain text (shared strings), matching the existing "numbers stored as
# text" convention already used in this sheet - NOT as numeric cells.
#
# Writing a string like "68" straight into .Value gets auto-coerced to a
# number by Excel's type inference, and forcing text via NumberFormat="@"
# (or a leading apostrophe) stamps a quote-prefix/text style onto the
# cell. Instead we park the literal text in a scratch cell via a formula
# (whose result is a string), copy it, and Paste Special "Values only"
# into the target - this carries over the text data type without
# touching any cell formatting/styles.
$newScores = @{
    1 = "68"
    2 = "85"
    3 = "56"
    4 = "83"
    5 = "66"
    6 = "73"
    7 = "67"
    8 = "59"
}

$scratch = $ws.Range("Z1")

foreach ($row in 1..8) {
    $scratch.Formula = "=""" + $newScores[$row] + """"
    $scratch.Copy()
    $ws.Range("C" + $row).PasteSpecial(-4163)
}

$scratch.Clear()
